$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.597787976264954
$ws.Range("B1").Value = 4.39522647857666
$ws.Range("C1").Value = 2.977920770645142
$ws.Range("D1").Value = 1.143656730651855
$ws.Range("E1").Value = 0.9038651585578918
